$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 needs to become a filled data row, matching the style/format of the
# other data rows (e.g. row 10). Copy formatting from row 10's B:E cells onto
# row 16's B:E cells first, then set the actual values.
$ws.Range("B10:E10").Copy() | Out-Null
$ws.Range("B16:E16").PasteSpecial(-4122) | Out-Null

$ws.Range("B16").Value = "TESTE EDSON"
$ws.Range("C16").Value = "530d170cb3487d322dbc4f7c8d853417"
$ws.Range("D16").Value = (Get-Date -Year 2022 -Month 9 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E16").Value = 8
